# chore: update Sheets via scheduled runner
# Refresh currentAveragePrice / LevePrice / LeveProfit columns (H:N) for the
# leves whose market-board prices moved since the last scheduled pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 100
$ws.Range("H100").Value = 1775.9286
$ws.Range("I100").Value = 1904.5
$ws.Range("J100").Value = 1454.5
$ws.Range("K100").Value = 1904.5
$ws.Range("L100").Value = 1454.5
$ws.Range("M100").Value = -1363.5
$ws.Range("N100").Value = -2536.5

# Row 117
$ws.Range("H117").Value = 48307.332
$ws.Range("J117").Value = 48307.332
$ws.Range("L117").Value = 48307.332
$ws.Range("N117").Value = -57485.332

# Row 131
$ws.Range("H131").Value = 2714.7144
$ws.Range("J131").Value = 3135.7144
$ws.Range("L131").Value = 9407.143199999999
$ws.Range("N131").Value = -19487.1432

# Row 138
$ws.Range("H138").Value = 1547.07
$ws.Range("I138").Value = 941.58826
$ws.Range("J138").Value = 1671.0844
$ws.Range("K138").Value = 2824.76478
$ws.Range("L138").Value = 5013.2532
$ws.Range("M138").Value = 2315.23522
$ws.Range("N138").Value = -15293.2532

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2644.476
$ws.Range("I61").Value = 1522.88
$ws.Range("J61").Value = 4293.8823
$ws.Range("K61").Value = 1522.88
$ws.Range("L61").Value = 4293.8823
$ws.Range("M61").Value = -1310.88
$ws.Range("N61").Value = -4717.8823

# Row 104
$ws.Range("H104").Value = 40732.668
$ws.Range("J104").Value = 40732.668
$ws.Range("L104").Value = 40732.668
$ws.Range("N104").Value = -47720.668

# Row 105
$ws.Range("H105").Value = 47967
$ws.Range("J105").Value = 47967
$ws.Range("L105").Value = 47967
$ws.Range("N105").Value = -54955

# Row 110
$ws.Range("H110").Value = 1944.4375
$ws.Range("I110").Value = 1944.4375
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1944.4375
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 100.5625
$ws.Range("N110").ClearContents()

# Row 117
$ws.Range("H117").Value = 48361
$ws.Range("J117").Value = 48361
$ws.Range("L117").Value = 48361
$ws.Range("N117").Value = -57539

# Row 132
$ws.Range("H132").Value = 19232314
$ws.Range("I132").Value = 26316824
$ws.Range("J132").Value = 2928
$ws.Range("K132").Value = 78950472
$ws.Range("L132").Value = 8784
$ws.Range("M132").Value = -78947942
$ws.Range("N132").Value = -13844

# Row 136
$ws.Range("H136").Value = 2644.476
$ws.Range("I136").Value = 1522.88
$ws.Range("J136").Value = 4293.8823
$ws.Range("K136").Value = 4568.64
$ws.Range("L136").Value = 12881.6469
$ws.Range("M136").Value = -2018.64
$ws.Range("N136").Value = -17981.6469

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 602.9091
$ws.Range("I94").Value = 614.0526
$ws.Range("J94").Value = 532.3333
$ws.Range("K94").Value = 614.0526
$ws.Range("L94").Value = 532.3333
$ws.Range("M94").Value = -163.0526
$ws.Range("N94").Value = -1434.3333

# Row 134
$ws.Range("H134").Value = 3534.7126
$ws.Range("I134").Value = 1712.4615
$ws.Range("K134").Value = 5137.3845
$ws.Range("M134").Value = -2602.3845

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 169139.28
$ws.Range("I31").Value = 1963.7778
$ws.Range("J31").Value = 214052.11
$ws.Range("K31").Value = 1963.7778
$ws.Range("L31").Value = 214052.11
$ws.Range("M31").Value = -1668.7778
$ws.Range("N31").Value = -214642.11

# Row 34
$ws.Range("H34").Value = 169139.28
$ws.Range("I34").Value = 1963.7778
$ws.Range("J34").Value = 214052.11
$ws.Range("K34").Value = 1963.7778
$ws.Range("L34").Value = 214052.11
$ws.Range("M34").Value = -1761.7778
$ws.Range("N34").Value = -214456.11

# Row 99
$ws.Range("H99").Value = 1949.5714
$ws.Range("I99").Value = 2049.5715
$ws.Range("J99").Value = 1899.5714
$ws.Range("K99").Value = 2049.5715
$ws.Range("L99").Value = 1899.5714
$ws.Range("M99").Value = -551.5715
$ws.Range("N99").Value = -4895.5714

# Row 109
$ws.Range("H109").Value = 28715.223
$ws.Range("J109").Value = 28715.223
$ws.Range("L109").Value = 28715.223
$ws.Range("N109").Value = -30795.223

# Row 111
$ws.Range("H111").Value = 47276.332
$ws.Range("J111").Value = 47276.332
$ws.Range("L111").Value = 47276.332
$ws.Range("N111").Value = -55456.332

# Row 116
$ws.Range("H116").Value = 47823.5
$ws.Range("J116").Value = 47823.5
$ws.Range("L116").Value = 47823.5
$ws.Range("N116").Value = -57001.5

# Row 126
$ws.Range("H126").Value = 1949.5714
$ws.Range("I126").Value = 2049.5715
$ws.Range("J126").Value = 1899.5714
$ws.Range("K126").Value = 6148.7145
$ws.Range("L126").Value = 5698.7142
$ws.Range("M126").Value = -3678.7145
$ws.Range("N126").Value = -10638.7142

# Row 141
$ws.Range("H141").Value = 12281.286
$ws.Range("J141").Value = 12281.286
$ws.Range("L141").Value = 12281.286
$ws.Range("N141").Value = -22641.286

$ws = $wb.Worksheets.Item("GSM")
# Row 21
$ws.Range("H21").Value = 5991.4287
$ws.Range("I21").Value = 1000
$ws.Range("J21").Value = 6823.3335
$ws.Range("K21").Value = 1000
$ws.Range("L21").Value = 6823.3335
$ws.Range("M21").Value = -827
$ws.Range("N21").Value = -7169.3335

# Row 26
$ws.Range("H26").Value = 22330
$ws.Range("J26").Value = 22330
$ws.Range("L26").Value = 22330
$ws.Range("N26").Value = -22890

# Row 30
$ws.Range("H30").Value = 5991.4287
$ws.Range("I30").Value = 1000
$ws.Range("J30").Value = 6823.3335
$ws.Range("K30").Value = 1000
$ws.Range("L30").Value = 6823.3335
$ws.Range("M30").Value = -895
$ws.Range("N30").Value = -7033.3335

# Row 50
$ws.Range("H50").Value = 22330
$ws.Range("J50").Value = 22330
$ws.Range("L50").Value = 22330
$ws.Range("N50").Value = -23326

# Row 97
$ws.Range("H97").Value = 5911.778
$ws.Range("I97").Value = 4742.143
$ws.Range("J97").Value = 10005.5
$ws.Range("K97").Value = 4742.143
$ws.Range("L97").Value = 10005.5
$ws.Range("M97").Value = -4246.143
$ws.Range("N97").Value = -10997.5

# Row 104
$ws.Range("H104").Value = 44961.5
$ws.Range("J104").Value = 44961.5
$ws.Range("L104").Value = 44961.5
$ws.Range("N104").Value = -51949.5

# Row 105
$ws.Range("H105").Value = 43998
$ws.Range("J105").Value = 43998
$ws.Range("L105").Value = 43998
$ws.Range("N105").Value = -50986

# Row 120
$ws.Range("H120").Value = 33531.332
$ws.Range("J120").Value = 33531.332
$ws.Range("L120").Value = 33531.332
$ws.Range("N120").Value = -43207.332

$ws = $wb.Worksheets.Item("LTW")
# Row 97
$ws.Range("H97").Value = 35000
$ws.Range("J97").Value = 35000
$ws.Range("L97").Value = 35000
$ws.Range("N97").Value = -36982

# Row 110
$ws.Range("H110").Value = 45581
$ws.Range("J110").Value = 45581
$ws.Range("L110").Value = 45581
$ws.Range("N110").Value = -53761

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 1497.625
$ws.Range("I96").Value = 998
$ws.Range("J96").Value = 1569
$ws.Range("K96").Value = 998
$ws.Range("L96").Value = 1569
$ws.Range("M96").Value = 375
$ws.Range("N96").Value = -4315
